$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# Remove the stored automatic-reply text from column E (row 2)
$ws.Range("E2").ClearContents()

# Update the timestamp of the log entry
$ws.Range("F2").Value = "2025-07-27 18:30:31"

# Update the status flags for this row
$ws.Range("G2").Value = "Nee"
$ws.Range("H2").Value = "Ja"
$ws.Range("I2").Value = "Nee"
